# Applies the update described by the commit:
#  - Swap the F:V (match details/odds/url) content between four row pairs
#    (the match order for same-date fixtures was corrected)
#  - Append two new fixtures as rows 153 and 154

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($sheet, [int]$row1, [int]$row2, [string]$colStart, [string]$colEnd)

    $r1 = $sheet.Range("$colStart$row1`:$colEnd$row1")
    $r2 = $sheet.Range("$colStart$row2`:$colEnd$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# --- Swap the four mis-ordered fixture pairs (columns F..V only; A..E —
#     index/country/league/season/date — stay put) ---
Swap-RowRange $ws 32  33  "F" "V"
Swap-RowRange $ws 37  38  "F" "V"
Swap-RowRange $ws 125 126 "F" "V"
Swap-RowRange $ws 133 134 "F" "V"

function Set-TextCell {
    param($sheet, [string]$addr, [string]$text)

    # Force a genuinely numeric-looking string (e.g. "2023") to be stored
    # as text rather than being auto-coerced to a number by the COM
    # value-setter, then drop the temporary text format so the cell is
    # left with the sheet's default style (matches the rest of column D).
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.ClearFormats()
}

function Add-Fixture {
    param(
        $sheet,
        [int]$rowIndex,
        [int]$idx,
        [string]$country,
        [string]$league,
        [string]$season,
        [double]$date,
        [string]$home,
        [int]$homeGoals,
        [string]$away,
        [int]$awayGoals,
        [double]$odds1Open, [string]$odds1OpenDate, [double]$odds1Close, [string]$odds1CloseDate,
        [double]$oddsXOpen, [string]$oddsXOpenDate, [double]$oddsXClose, [string]$oddsXCloseDate,
        [double]$odds2Open, [string]$odds2OpenDate, [double]$odds2Close, [string]$odds2CloseDate,
        [string]$url
    )

    $prevRow = $rowIndex - 1

    # Carry the row-level styles (bold/border index cell in col A, the
    # date/time number format in col E) from the previous (last existing)
    # data row onto the freshly appended one.
    $sheet.Range("A$prevRow`:V$prevRow").Copy()
    $sheet.Range("A$rowIndex`:V$rowIndex").PasteSpecial(-4122) # xlPasteFormats

    $sheet.Range("A$rowIndex").Value2 = $idx
    $sheet.Range("B$rowIndex").Value2 = $country
    $sheet.Range("C$rowIndex").Value2 = $league
    Set-TextCell $sheet "D$rowIndex" $season
    $sheet.Range("E$rowIndex").Value2 = $date
    $sheet.Range("F$rowIndex").Value2 = $home
    $sheet.Range("G$rowIndex").Value2 = $homeGoals
    $sheet.Range("H$rowIndex").Value2 = $away
    $sheet.Range("I$rowIndex").Value2 = $awayGoals
    $sheet.Range("J$rowIndex").Value2 = $odds1Open
    $sheet.Range("K$rowIndex").Value2 = $odds1OpenDate
    $sheet.Range("L$rowIndex").Value2 = $odds1Close
    $sheet.Range("M$rowIndex").Value2 = $odds1CloseDate
    $sheet.Range("N$rowIndex").Value2 = $oddsXOpen
    $sheet.Range("O$rowIndex").Value2 = $oddsXOpenDate
    $sheet.Range("P$rowIndex").Value2 = $oddsXClose
    $sheet.Range("Q$rowIndex").Value2 = $oddsXCloseDate
    $sheet.Range("R$rowIndex").Value2 = $odds2Open
    $sheet.Range("S$rowIndex").Value2 = $odds2OpenDate
    $sheet.Range("T$rowIndex").Value2 = $odds2Close
    $sheet.Range("U$rowIndex").Value2 = $odds2CloseDate
    $sheet.Range("V$rowIndex").Value2 = $url
}

Add-Fixture $ws 153 152 "argentina" "copa-de-la-liga-profesional" "2023" 45230.95833333334 `
    "Tigre" 1 "Godoy Cruz" 0 `
    2.04 "27/10/2023 02:12" 2.66 "31/10/2023 22:59" `
    3.29 "27/10/2023 02:12" 3.01 "31/10/2023 22:59" `
    3.99 "27/10/2023 02:12" 3.08 "31/10/2023 22:57" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/tigre-godoy-cruz/MDs09Pyp/"

Add-Fixture $ws 154 153 "argentina" "copa-de-la-liga-profesional" "2023" 45231.04166666666 `
    "Central Cordoba" 2 "Union de Santa Fe" 0 `
    2.45 "25/10/2023 23:42" 2.99 "01/11/2023 00:57" `
    3.04 "25/10/2023 23:42" 2.72 "01/11/2023 00:57" `
    3.29 "25/10/2023 23:42" 3.02 "01/11/2023 00:57" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/central-cordoba-santiago-del-estero-union-de-santa-fe/QsIaUQD9/"
